$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: numeric value 0
$ws.Range("B1").Value = 0

# Bold font, thin box border, center/top alignment on B1
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B1").VerticalAlignment = -4160    # xlTop
$ws.Range("B1").Borders.LineStyle = 1        # xlContinuous
$ws.Range("B1").Borders.Weight = 2           # xlThin

# A2: same value + identical style as B1 (copy formats so the same cellXf is reused)
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# B2: text value, default (unstyled) cell
$ws.Range("B2").Value = "disconnected_elements"
